$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the target paragraph: "Future implementation: Maybe allow
# for different structures" (currently one run, with a stray
# "_GoBack" bookmark splitting "...structu" / "res" as an artifact of
# the author's last edit position).
# ------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt.StartsWith("Future implementation:")) {
        $targetIndex = $i
        break
    }
}

$boldLead = "Future implementation:"
$suggestionText = "Suggestion: Ask company to specify column names for their senders, receivers, transactions"

$p = $d.Paragraphs($targetIndex)
$pStart = $p.Range.Start

# ------------------------------------------------------------------
# 1) Make the "Future implementation:" lead-in its own bold run.
# ------------------------------------------------------------------
$boldRange = $d.Range($pStart, $pStart + $boldLead.Length)
$boldRange.Font.Bold = 1

# ------------------------------------------------------------------
# 2) Insert a brand-new list paragraph right after it, one level
#    deeper (ilvl 2, same numId) carrying the "Suggestion: ..." text,
#    with the "_GoBack" bookmark re-anchored at the end of it (where
#    the author's cursor ended up after typing the new line).
#
#    A sacrificial trailing "Z" is appended before placing the
#    bookmark: collapsing a freshly built Range exactly on
#    "paragraph end - 1" resolves to the wrong spot in this host, but
#    anchoring one character earlier (still inside real text) and
#    trimming the sacrificial character afterwards keeps the bookmark
#    glued to the correct place, immediately after the real text.
# ------------------------------------------------------------------
$d1 = $word.ActiveDocument
$srcPara = $d1.Paragraphs($targetIndex)
$srcPara.Range.InsertParagraphAfter()

$d2 = $word.ActiveDocument
$newPara = $d2.Paragraphs($targetIndex + 1)
$newPara.Range.ListFormat.ListLevelNumber = 3
$newPara.Range.Text = $suggestionText + "Z"

$d3 = $word.ActiveDocument
$newParaRange = $d3.Paragraphs($targetIndex + 1).Range
$bmPos = $newParaRange.End - 2
$bmRange = $d3.Range($bmPos, $bmPos)
$d3.Bookmarks.Add("_GoBack", $bmRange)

$d4 = $word.ActiveDocument
$newParaRange2 = $d4.Paragraphs($targetIndex + 1).Range
$zPos = $newParaRange2.End - 2
$zRange = $d4.Range($zPos, $zPos + 1)
$zRange.Text = ""

# ------------------------------------------------------------------
# 3) Back in the original paragraph, the old bookmark used to sit
#    between "...structu" and "res", forcing two runs. Now that it
#    has moved to the new paragraph, rewrite the (non-bold) tail of
#    the sentence from scratch so it collapses back into a single
#    run reading "...structures".
# ------------------------------------------------------------------
$d5 = $word.ActiveDocument
$origPara = $d5.Paragraphs($targetIndex)
$tailStart = $origPara.Range.Start + $boldLead.Length
$tailEnd = $origPara.Range.End - 1
$tailRange = $d5.Range($tailStart, $tailEnd)
$tailRange.Text = ""

$d6 = $word.ActiveDocument
$origPara2 = $d6.Paragraphs($targetIndex)
$insPos = $origPara2.Range.Start + $boldLead.Length
$insRange = $d6.Range($insPos, $insPos)
$insRange.Text = " Maybe allow for different structures"
$insRange.Font.Bold = 0

# ------------------------------------------------------------------
# 4) The third-level bullet ("ilvl 2") of the list used by this
#    section is no longer a "tentative" preview level once it is
#    actually used in the body text -- clear the tentative flag on
#    its numbering definition.
# ------------------------------------------------------------------
$d7 = $word.ActiveDocument
$listTemplates = $d7.ListFormat.ListTemplate
$targetPara = $d7.Paragraphs($targetIndex + 1)
$targetPara.Range.ListFormat.ListTemplate.ListLevels(3).Tentative = 0

Write-Host "Edit complete"
